$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.111.45'
$ws.Range('E2').Value = '  +5.11%  '
$ws.Range('D3').Value = '2.439.97'
$ws.Range('E3').Value = '  +5.48%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''565.08'
$ws.Range('E5').Value = '  +4.15%  '
$ws.Range('D6').Value = '''140.99'
$ws.Range('E6').Value = '  +8.40%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '''0.587'
$ws.Range('E8').Value = '  +2.57%  '
$ws.Range('D9').Value = '2.438.15'
$ws.Range('E9').Value = '  +5.51%  '
$ws.Range('E10').Value = '  +3.87%  '
$ws.Range('D11').Value = '''5.76'
$ws.Range('E11').Value = '  +3.72%  '
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('D13').Value = '''0.351'
$ws.Range('E13').Value = '  +5.75%  '
$ws.Range('D14').Value = '''26.37'
$ws.Range('E14').Value = '  +12.99%  '
$ws.Range('D15').Value = '2.873.41'
$ws.Range('E15').Value = '  +5.49%  '
$ws.Range('D16').Value = '62.975.25'
$ws.Range('E16').Value = '  +5.01%  '
$ws.Range('D17').Value = '''0.0000143'
$ws.Range('E17').Value = '  +7.99%  '
$ws.Range('D18').Value = '2.439.96'
$ws.Range('E18').Value = '  +5.95%  '
$ws.Range('D19').Value = '''11.29'
$ws.Range('E19').Value = '  +7.71%  '
$ws.Range('D20').Value = '''340.79'
$ws.Range('E20').Value = '  +9.28%  '
$ws.Range('D21').Value = '''4.23'
$ws.Range('E21').Value = '  +4.13%  '
$ws.Range('D22').Value = '''6.82'
$ws.Range('E22').Value = '  +4.35%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = '''5.65'
$ws.Range('E24').Value = '  -0.59%  '
$ws.Range('D25').Value = '''65.55'
$ws.Range('E25').Value = '  +3.42%  '
$ws.Range('D26').Value = '''0.174'
$ws.Range('E26').Value = '  +2.90%  '
$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').Value = '''1.55'
$ws.Range('E27').Value = '  +14.70%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = '''1.00'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').Value = '''8.24'
$ws.Range('E29').Value = '  +6.34%  '
$ws.Range('E30').Value = '  +14.14%  '
$ws.Range('D31').Value = '0.0₃0796'
$ws.Range('E31').Value = '  +9.78%  '
$ws.Range('E32').Value = '  +4.98%  '
$ws.Range('E33').Value = '  +12.16%  '
$ws.Range('D34').Value = '''174.28'
$ws.Range('E34').Value = '  +1.55%  '
$ws.Range('E35').Value = '  +9.40%  '
$ws.Range('D36').Value = '''0.399'
$ws.Range('E36').Value = '  +5.33%  '
$ws.Range('D37').Value = '''380.40'
$ws.Range('E37').Value = '  +19.64%  '
$ws.Range('D38').Value = '''18.67'
$ws.Range('E38').Value = '  +5.61%  '
$ws.Range('D39').Value = '''4.50'
$ws.Range('E39').Value = '  +11.99%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').Value = '''1.72'
$ws.Range('E42').Value = '  +13.58%  '
$ws.Range('D43').Value = '''39.97'
$ws.Range('E43').Value = '  +6.61%  '
$ws.Range('D44').Value = '''145.43'
$ws.Range('E44').Value = '  +6.60%  '
$ws.Range('E45').Value = '  +6.98%  '
$ws.Range('D46').Value = '''20.70'
$ws.Range('E46').Value = '  +10.41%  '
$ws.Range('D47').Value = '''0.595'
$ws.Range('E47').Value = '  +4.67%  '
$ws.Range('E48').Value = '  +6.14%  '
$ws.Range('D49').Value = '''0.0949'
$ws.Range('E49').Value = '  +0.54%  '
$ws.Range('E50').Value = '  +5.07%  '
$ws.Range('D51').Value = '''17.91'
$ws.Range('E51').Value = '  +7.11%  '
